# "#5: property building done" — append the data row (row 2) under the
# existing header-ish row (row 1) on the single worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A2: numeric 28, formatted like the row-1 cells (bordered/centered style).
$ws.Range("A2").Value = 28
$ws.Range("B1").Copy()
$ws.Range("A2").PasteSpecial(-4122)   # xlPasteFormats

# B2:G2: same values as B1:G1 (reuses the existing shared strings / numbers).
$ws.Range("B2").Value = $ws.Range("B1").Value2
$ws.Range("C2").Value = $ws.Range("C1").Value2
$ws.Range("D2").Value = $ws.Range("D1").Value2
$ws.Range("E2").Value = $ws.Range("E1").Value2
$ws.Range("F2").Value = $ws.Range("F1").Value2
$ws.Range("G2").Value = $ws.Range("G1").Value2

# Give B2:G2 their own (plain/default-looking) style distinct from row 1's.
$ws.Range("B2:G2").HorizontalAlignment = 1   # xlGeneral — default, but mints its own style record
